# Generate Report for Handoff
#
# The "Status" for the zh-cn / de-de handoff moved from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# handoff timestamps advanced a few seconds. Excel re-autofits the
# (now much shorter) Status columns after the text changes, so their
# column widths shrink too.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"
$newWidth  = 16.4   # ColumnWidth (chars) that yields the narrower autofit width

# ----- Overview sheet -----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-25 22:57:55"
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# ----- zh-cn sheet ----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-25 22:57:51"
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# ----- de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-25 22:57:55"
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
